# #200: Excel date format config (#294)
# Config options to control Excel date/time output format
#
# The "col_time" / "col_time_n" columns (P:Q) were previously written out as
# text strings ("00:00:00", "07:07:07"). With the new date/time formatting
# config they are now written as real numeric Excel time-of-day serial
# values, displayed with a "hh:mm:ss" number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> col_time / col_time_n : "00:00:00" (text) becomes 0 (numeric time)
$ws.Range("P2:Q2").NumberFormat = "hh:mm:ss"
$ws.Range("P2:Q2").Value = 0

# Row 3 -> col_time / col_time_n : "07:07:07" (text) becomes the matching
# Excel time-of-day fraction (same fractional part already used by the
# col_datetime / col_datetime_n columns on this row).
$ws.Range("P3:Q3").NumberFormat = "hh:mm:ss"
$ws.Range("P3:Q3").Value = 0.2966087962995516

# Row 4 -> col_time only (col_time_n stays blank on this row, as before)
$ws.Range("P4").NumberFormat = "hh:mm:ss"
$ws.Range("P4").Value = 0.2966087962995516

Write-Host "Applied Excel date/time format config changes to col_time/col_time_n"
